# Dagbaekur Ivan Gusti Saevar uddfaerdar
# Updates weekly timesheet entries (Vika 5 / Vika 6) and the summary formula.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# --- Vika 4 (row 17): þri logged 90 min extra research time ---
$ws.Range("E17").Value = 120

# --- Vika 5 (rows 36-41): fill in the previously-empty week ---
$ws.Range("C36").Value = 60
$ws.Range("D36").Value = 60

$ws.Range("D37").Value = 60

$ws.Range("D38").Value = 120
$ws.Range("E38").Value = 60

$ws.Range("C39").Value = 120
$ws.Range("D39").Value = 280
$ws.Range("E39").Value = 90

# --- Vika 6 (rows 46-51): fill in the previously-empty week ---
$ws.Range("C46").Value = 90

$ws.Range("D49").Value = 120

$ws.Range("C50").Value = 30
$ws.Range("D50").Value = 180
$ws.Range("E50").Value = 120

$ws.Range("D51").Value = 60
$ws.Range("E51").Value = 240

# --- Samantekt (row 62): total now sums the weekly subtotals directly,
#     and a new helper column converts the grand total from minutes to hours ---
$ws.Range("D62").Formula = "=SUM(J52,J42,J32,J22,J12)"
$ws.Range("F62").Formula = "=D62/60"

# --- view state: scroll down to week 6/7 and select the newly-filled block ---
$win = $excel.ActiveWindow
$win.ScrollRow = 52
$win.ScrollColumn = 1
[void]$ws.Range("C36:I41").Select()
